# Refresh the market-price-derived columns (H:N) for the leve rows the
# scheduled price-sync run touched, across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 7093.3
$ws.Range("I62").Value = 7376.5454
$ws.Range("J62").Value = 5758
$ws.Range("K62").Value = 7376.5454
$ws.Range("L62").Value = 5758
$ws.Range("M62").Value = -6752.5454
$ws.Range("N62").Value = -7006
# Row 65
$ws.Range("H65").Value = 7093.3
$ws.Range("I65").Value = 7376.5454
$ws.Range("J65").Value = 5758
$ws.Range("K65").Value = 36882.727
$ws.Range("L65").Value = 28790
$ws.Range("M65").Value = -33762.727
$ws.Range("N65").Value = -35030
# Row 116
$ws.Range("H116").Value = 33792
$ws.Range("I116").Value = 55332
$ws.Range("J116").Value = 3020.5715
$ws.Range("K116").Value = 55332
$ws.Range("L116").Value = 3020.5715
$ws.Range("M116").Value = -51890
$ws.Range("N116").Value = -9904.5715
# Row 117
$ws.Range("H117").Value = 31471
$ws.Range("J117").Value = 31471
$ws.Range("L117").Value = 31471
$ws.Range("N117").Value = -40649
# Row 127
$ws.Range("H127").Value = 931.0769
$ws.Range("I127").Value = 487.44446
$ws.Range("J127").Value = 1165.9412
$ws.Range("K127").Value = 1462.33338
$ws.Range("L127").Value = 3497.8236
$ws.Range("M127").Value = 3497.66662
$ws.Range("N127").Value = -13417.8236
# Row 132
$ws.Range("H132").Value = 1862.0886
$ws.Range("I132").Value = 1024.3572
$ws.Range("J132").Value = 8377.777
$ws.Range("K132").Value = 3073.0716
$ws.Range("L132").Value = 25133.331
$ws.Range("M132").Value = -543.0715999999998
$ws.Range("N132").Value = -30193.331
# Row 137
$ws.Range("H137").Value = 2627.9473
$ws.Range("I137").Value = 2702.2307
$ws.Range("K137").Value = 8106.6921
$ws.Range("M137").Value = -5556.6921

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1237.2174
$ws.Range("I45").Value = 978.2308
$ws.Range("J45").Value = 1573.9
$ws.Range("K45").Value = 978.2308
$ws.Range("L45").Value = 1573.9
$ws.Range("M45").Value = -601.2308
$ws.Range("N45").Value = -2327.9

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 759.53845
$ws.Range("I94").Value = 644.069
$ws.Range("J94").Value = 1094.4
$ws.Range("K94").Value = 644.069
$ws.Range("L94").Value = 1094.4
$ws.Range("M94").Value = -193.069
$ws.Range("N94").Value = -1996.4
# Row 99
$ws.Range("H99").Value = 3234.3635
$ws.Range("I99").Value = 4720.6924
$ws.Range("J99").Value = 1087.4445
$ws.Range("K99").Value = 4720.6924
$ws.Range("L99").Value = 1087.4445
$ws.Range("M99").Value = -3222.6924
$ws.Range("N99").Value = -4083.4445

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 916.9375
$ws.Range("I16").Value = 898.9231
$ws.Range("J16").Value = 995
$ws.Range("K16").Value = 898.9231
$ws.Range("L16").Value = 995
$ws.Range("M16").Value = -611.9231
$ws.Range("N16").Value = -1569
# Row 31
$ws.Range("H31").Value = 2010.0127
$ws.Range("I31").Value = 1422.4445
$ws.Range("J31").Value = 3279.16
$ws.Range("K31").Value = 1422.4445
$ws.Range("L31").Value = 3279.16
$ws.Range("M31").Value = -1127.4445
$ws.Range("N31").Value = -3869.16
# Row 34
$ws.Range("H34").Value = 2010.0127
$ws.Range("I34").Value = 1422.4445
$ws.Range("J34").Value = 3279.16
$ws.Range("K34").Value = 1422.4445
$ws.Range("L34").Value = 3279.16
$ws.Range("M34").Value = -1220.4445
$ws.Range("N34").Value = -3683.16
# Row 60
$ws.Range("H60").Value = 10950
$ws.Range("I60").Value = 9996.5
$ws.Range("J60").Value = 11050.368
$ws.Range("K60").Value = 9996.5
$ws.Range("L60").Value = 11050.368
$ws.Range("M60").Value = -9485.5
$ws.Range("N60").Value = -12072.368
# Row 105
$ws.Range("H105").Value = 1473.5
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 1757.6
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 1757.6
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -5251.6
# Row 107
$ws.Range("H107").Value = 166.60715
$ws.Range("I107").Value = 144.73914
$ws.Range("J107").Value = 267.2
$ws.Range("K107").Value = 144.73914
$ws.Range("L107").Value = 267.2
$ws.Range("M107").Value = 1775.26086
$ws.Range("N107").Value = -4107.2
# Row 113
$ws.Range("H113").Value = 916.9375
$ws.Range("I113").Value = 898.9231
$ws.Range("J113").Value = 995
$ws.Range("K113").Value = 898.9231
$ws.Range("L113").Value = 995
$ws.Range("M113").Value = 1271.0769
$ws.Range("N113").Value = -5335
# Row 134
$ws.Range("H134").Value = 1392.0212
$ws.Range("I134").Value = 919.45715
$ws.Range("K134").Value = 2758.37145
$ws.Range("M134").Value = -223.3714499999996

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 4310.1816
$ws.Range("I87").Value = 2844.5715
$ws.Range("K87").Value = 8533.7145
$ws.Range("M87").Value = -7285.7145
# Row 90
$ws.Range("H90").Value = 4310.1816
$ws.Range("I90").Value = 2844.5715
$ws.Range("K90").Value = 25601.1435
$ws.Range("M90").Value = -19361.1435
# Row 92
$ws.Range("H92").Value = 33333538
$ws.Range("I92").Value = 38461730
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 115385190
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = -115383942
$ws.Range("N92").Value = -3396
# Row 103
$ws.Range("H103").Value = 4157
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 5209.3335
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 15628.0005
$ws.Range("M103").Value = -2121
$ws.Range("N103").Value = -17386.0005
# Row 113
$ws.Range("H113").Value = 9091454
$ws.Range("I113").Value = 13158427
$ws.Range("K113").Value = 39475281
$ws.Range("M113").Value = -39473111

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 126
$ws.Range("H126").Value = 3310.1924
$ws.Range("I126").Value = 3195
$ws.Range("J126").Value = 3361.389
$ws.Range("K126").Value = 9585
$ws.Range("L126").Value = 10084.167
$ws.Range("M126").Value = -7115
$ws.Range("N126").Value = -15024.167

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1609.375
$ws.Range("I61").Value = 1696.4286
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 1696.4286
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -1494.4286
$ws.Range("N61").Value = -1404
# Row 113
$ws.Range("H113").Value = 1609.375
$ws.Range("I113").Value = 1696.4286
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1696.4286
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 473.5714
$ws.Range("N113").Value = -5340
# Row 122
$ws.Range("H122").Value = 2139.28
$ws.Range("I122").Value = 2282.611
$ws.Range("J122").Value = 1770.7142
$ws.Range("K122").Value = 6847.833
$ws.Range("L122").Value = 5312.142599999999
$ws.Range("M122").Value = -4397.833
$ws.Range("N122").Value = -10212.1426
